# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the Hades_Profits (FFXIV leve profit) workbook
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# -- ALC row 8 (On the Drip) --
$ws.Range("H8").Value = 2411.8
$ws.Range("I8").Value = 29.5
$ws.Range("J8").Value = 4000
$ws.Range("K8").Value = 88.5
$ws.Range("L8").Value = 12000
$ws.Range("M8").Value = 50.5
$ws.Range("N8").Value = -12278

# -- ALC row 40 (Stuck in the Moment) --
$ws.Range("H40").Value = 1000
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

# -- ALC row 100 (Asking for a Friend) --
$ws.Range("H100").Value = 2780
$ws.Range("I100").Value = 2725
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2725
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -2184
$ws.Range("N100").Value = -4082

# -- ALC row 103 (Let Loose the Juice) --
$ws.Range("H103").Value = 1000
$ws.Range("I103").Value = 733.3333
$ws.Range("J103").Value = 1400
$ws.Range("K103").Value = 2199.9999
$ws.Range("L103").Value = 4200
$ws.Range("M103").Value = -1613.9999
$ws.Range("N103").Value = -5372

$ws = $wb.Worksheets.Item("ARM")
# -- ARM row 32 (Ingot We Trust) --
$ws.Range("H32").Value = 4646220.5
$ws.Range("I32").Value = 5145970
$ws.Range("K32").Value = 5145970
$ws.Range("M32").Value = -5145683

# -- ARM row 74 (As the Bolt Flies) --
$ws.Range("H74").Value = 12001299
$ws.Range("I74").Value = 22818806
$ws.Range("J74").Value = 102040.4
$ws.Range("K74").Value = 22818806
$ws.Range("L74").Value = 102040.4
$ws.Range("M74").Value = -22817932
$ws.Range("N74").Value = -103788.4

# -- ARM row 77 (Heavy Metal Banned (L)) --
$ws.Range("H77").Value = 12001299
$ws.Range("I77").Value = 22818806
$ws.Range("J77").Value = 102040.4
$ws.Range("K77").Value = 114094030
$ws.Range("L77").Value = 510202
$ws.Range("M77").Value = -114089662
$ws.Range("N77").Value = -518938

# -- ARM row 109 (A Head of Demand) --
$ws.Range("H109").Value = 28600
$ws.Range("J109").Value = 28600
$ws.Range("L109").Value = 28600
$ws.Range("N109").Value = -31374

# -- ARM row 112 (Wrapped Knuckles) --
$ws.Range("H112").Value = 16600
$ws.Range("J112").Value = 16600
$ws.Range("L112").Value = 16600
$ws.Range("N112").Value = -19554

$ws = $wb.Worksheets.Item("BSM")
# -- BSM row 12 (A Hit Job) --
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

# -- BSM row 105 (Ingot to Wing It) --
$ws.Range("H105").Value = 55557720
$ws.Range("I105").Value = 100002296
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 100002296
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -100000549
$ws.Range("N105").Value = -5494

$ws = $wb.Worksheets.Item("CRP")
# -- CRP row 58 (You Do the Heavy Lifting) --
$ws.Range("H58").Value = 30912398
$ws.Range("I58").Value = 37780708
$ws.Range("J58").Value = 5002.3335
$ws.Range("K58").Value = 37780708
$ws.Range("L58").Value = 5002.3335
$ws.Range("M58").Value = -37780505
$ws.Range("N58").Value = -5408.3335

# -- CRP row 132 (Hull Lotta Damage) --
$ws.Range("H132").Value = 77550.53999999999
$ws.Range("I132").Value = 650.63635
$ws.Range("J132").Value = 500500
$ws.Range("K132").Value = 1951.90905
$ws.Range("L132").Value = 1501500
$ws.Range("M132").Value = 578.09095
$ws.Range("N132").Value = -1506560

# -- CRP row 134 (Wood You Be Quiet) --
$ws.Range("H134").Value = 44826.72
$ws.Range("I134").Value = 1577.5
$ws.Range("J134").Value = 121714.22
$ws.Range("K134").Value = 4732.5
$ws.Range("L134").Value = 365142.66
$ws.Range("M134").Value = -2197.5
$ws.Range("N134").Value = -370212.66

# -- CRP row 136 (Turali Quality) --
$ws.Range("H136").Value = 30912398
$ws.Range("I136").Value = 37780708
$ws.Range("J136").Value = 5002.3335
$ws.Range("K136").Value = 113342124
$ws.Range("L136").Value = 15007.0005
$ws.Range("M136").Value = -113339574
$ws.Range("N136").Value = -20107.0005

$ws = $wb.Worksheets.Item("CUL")
# -- CUL row 68 (Such a Butter Face) --
$ws.Range("H68").Value = 972.13336
$ws.Range("I68").Value = 716.6667
$ws.Range("J68").Value = 1000.5185
$ws.Range("K68").Value = 2150.0001
$ws.Range("L68").Value = 3001.5555
$ws.Range("M68").Value = -1339.0001
$ws.Range("N68").Value = -4623.5555

# -- CUL row 71 (No Margarine of Error (L)) --
$ws.Range("H71").Value = 972.13336
$ws.Range("I71").Value = 716.6667
$ws.Range("J71").Value = 1000.5185
$ws.Range("K71").Value = 6450.0003
$ws.Range("L71").Value = 9004.666499999999
$ws.Range("M71").Value = -2394.0003
$ws.Range("N71").Value = -17116.6665

# -- CUL row 86 (Let's Not Get Sappy) --
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# -- CUL row 89 (Luxury Spillover (L)) --
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# -- CUL row 107 (Slippery Service) --
$ws.Range("H107").Value = 898.0893
$ws.Range("J107").Value = 1280.9032
$ws.Range("L107").Value = 3842.7096
$ws.Range("N107").Value = -7682.7096

# -- CUL row 122 (Salt of the North) --
$ws.Range("H122").Value = 620.2857
$ws.Range("I122").Value = 380.36365
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 3423.27285
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -973.2728500000003
$ws.Range("N122").Value = -18400

$ws = $wb.Worksheets.Item("GSM")
# -- GSM row 21 (Forever 21K) --
$ws.Range("H21").Value = 10416.667
$ws.Range("I21").Value = 2500
$ws.Range("J21").Value = 50000
$ws.Range("K21").Value = 2500
$ws.Range("L21").Value = 50000
$ws.Range("M21").Value = -2327
$ws.Range("N21").Value = -50346

# -- GSM row 30 (Dog Tags Are for Dogs) --
$ws.Range("H30").Value = 10416.667
$ws.Range("I30").Value = 2500
$ws.Range("J30").Value = 50000
$ws.Range("K30").Value = 2500
$ws.Range("L30").Value = 50000
$ws.Range("M30").Value = -2395
$ws.Range("N30").Value = -50210

# -- GSM row 31 (One and Only) --
$ws.Range("H31").Value = 1827.75
$ws.Range("I31").Value = 1827.75
$ws.Range("K31").Value = 1827.75
$ws.Range("M31").Value = -1535.75

# -- GSM row 37 (Dancing with the Stars) --
$ws.Range("H37").Value = 1827.75
$ws.Range("I37").Value = 1827.75
$ws.Range("K37").Value = 1827.75
$ws.Range("M37").Value = -1550.75

# -- GSM row 55 (If You've Got It, Flaunt It) --
$ws.Range("H55").Value = 1030
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

# -- GSM row 113 (Copious Crystal Cannons) --
$ws.Range("H113").Value = 1552.2
$ws.Range("I113").Value = 1552.2
$ws.Range("K113").Value = 1552.2
$ws.Range("M113").Value = 617.8

$ws = $wb.Worksheets.Item("LTW")
# -- LTW row 7 (Tan Before the Ban) --
$ws.Range("H7").Value = 3842
$ws.Range("I7").Value = 2987.3
$ws.Range("J7").Value = 5266.5
$ws.Range("K7").Value = 2987.3
$ws.Range("L7").Value = 5266.5
$ws.Range("M7").Value = -2875.3
$ws.Range("N7").Value = -5490.5

# -- LTW row 68 (You Could Say It's a Moving Target) --
$ws.Range("H68").Value = 1900
$ws.Range("I68").Value = 1800
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 1800
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -1051
$ws.Range("N68").Value = -3498

# -- LTW row 71 (They Call It Bloody Mary (L)) --
$ws.Range("H71").Value = 1900
$ws.Range("I71").Value = 1800
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -5256
$ws.Range("N71").Value = -17488

# -- LTW row 81 (I Need Your Glove Tonight) --
$ws.Range("H81").Value = 36620.668

# -- LTW row 84 (Halonic Drake Handlers (L)) --
$ws.Range("H84").Value = 36620.668

# -- LTW row 126 (Battered Books) --
$ws.Range("H126").Value = 3842
$ws.Range("I126").Value = 2987.3
$ws.Range("J126").Value = 5266.5
$ws.Range("K126").Value = 8961.900000000001
$ws.Range("L126").Value = 15799.5
$ws.Range("M126").Value = -6491.900000000001
$ws.Range("N126").Value = -20739.5

# -- LTW row 132 (Tenets of Tanning) --
$ws.Range("H132").Value = 44525.125
$ws.Range("I132").Value = 2417.4546
$ws.Range("J132").Value = 80154.69500000001
$ws.Range("K132").Value = 7252.3638
$ws.Range("L132").Value = 240464.085
$ws.Range("M132").Value = -4722.3638
$ws.Range("N132").Value = -245524.085

$ws = $wb.Worksheets.Item("WVR")
# -- WVR row 118 (Something in My Eye) --
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("M118").ClearContents()
$ws.Range("N118").ClearContents()

# -- WVR row 132 (Comfy Cabins) --
$ws.Range("H132").Value = 50356.22
$ws.Range("I132").Value = 39888.69
$ws.Range("J132").Value = 68499.92999999999
$ws.Range("K132").Value = 119666.07
$ws.Range("L132").Value = 205499.79
$ws.Range("M132").Value = -117136.07
$ws.Range("N132").Value = -210559.79

# -- WVR row 133 (Begin with the Basics) --
$ws.Range("H133").Value = 33476.668
$ws.Range("J133").Value = 33476.668
$ws.Range("L133").Value = 33476.668
$ws.Range("N133").Value = -43596.668
